# Update "想去人数" (F column) values on the "展览" and "全部类型" worksheets
# to reflect newly generated counts (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new F-column value (applies to both sheets; source values differ
# slightly between sheets but converge to these same targets)
$updates = @{
    2  = 308
    4  = 10311
    5  = 330
    6  = 941
    7  = 1279
    8  = 6997
    10 = 436
    11 = 198
    13 = 3176
    14 = 35
    15 = 311
    16 = 646
    17 = 122
    18 = 780
    20 = 60
    21 = 1618
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
